# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

# --- Rename "Requested quantity" headers so they're unambiguous per-sheet ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the existing sheets ---
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the bold/border/centered header style used on the other sheets
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data (dates in column A, numeric forecast + bounds in B:D)
$dates = @(45459.99999999999, 45466.99999999999, 45473.99999999999, 45515.99999999999, 45522.99999999999, 45529.99999999999, 45536.99999999999, 45543.99999999999, 45550.99999999999, 45557.99999999999, 45564.99999999999, 45571.99999999999, 45578.99999999999, 45585.99999999999)
$forecast = @(77, 72, 67, 37, 33, 28, 23, 18, 13, 8, 3, 0, 0, 0)
$lower = @(-13.77221160116944, -22.76492483936438, -22.2377034257627, -55.84364846529785, -61.25645478237598, -63.68676414329602, -71.44939217698936, -74.69099898240223, -78.97381462144986, -78.79382610771579, -86.12127103486331, -90.54484294770474, -96.4062641668819, -112.4253827685428)
$upper = @(168.241225172833, 164.2570623628673, 161.0651293038423, 126.5537571394522, 120.8134309468092, 114.2714014856075, 114.4703749648509, 104.7620506082488, 101.8387305614045, 95.60208109062806, 97.19077478763127, 90.04741514781004, 86.54714920706489, 81.58481452182052)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 1).Value2 = $dates[$i]
    $wsForecast.Cells.Item($row, 2).Value2 = $forecast[$i]
    $wsForecast.Cells.Item($row, 3).Value2 = $lower[$i]
    $wsForecast.Cells.Item($row, 4).Value2 = $upper[$i]
}

# Match the date-formatted style used for the date column on the other sheets
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A15").PasteSpecial(-4122)

$wsForecast.Range("A1").Select() | Out-Null
